$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Student Name -> Sahil Choudhary
$ws.Range("C3").Value = "Sahil Choudhary"

# Row 7 - Test case 1 (__init__, attributes set to parameter values)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number = 350`nclient_number = 350`nbalance = 350`ndate_created = (2024, 5, 10)`nminimum_balance = 30"
$ws.Range("G7").Value = "attributes are set"

# Row 8 - Test case 2 (__init__, minimum_balance has invalid type.)
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "account_number = 350`nclient_number = 350`nbalance = 350`ndate_created = (2024, 5, 10)`nminimum_balance = 'thirty'"
$ws.Range("G8").Value = "minimum_balance set to 100"

# Row 9 - Test case 3 (get_service_charges, balance greater than minimum balance)
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "account_number = 350`nclient_number = 350`nbalance = 350`ndate_created = (2024, 5, 10)`nminimum_balance = 30"
$ws.Range("G9").Value = "service_charge set to BASE_SERVICE_CHARGE "

# Row 10 - Test case 4 (get_service_charges, balance equal to minimum balance)
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "account_number = 350`nclient_number = 350`nbalance = 30`ndate_created = (2024, 5, 10)`nminimum_balance = 30"
$ws.Range("G10").Value = "service_charge set to BASE_SERVICE_CHARGE "

# Row 11 - Test case 5 (get_service_charges, balance less than minimum balance)
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = "account_number = 350`nclient_number = 350`nbalance = 15`ndate_created = (2024, 5, 10)`nminimum_balance = 30"
$ws.Range("G11").Value = "service_charge set to BASE_SERVICE_CHARGE * SERVICE_CHARGE_PREMIUM "

# Row 12 - Test case 6 (__str__, appropriate value returned based on attribute values.)
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "account_number = 350`nclient_number = 350`nbalance = 350`ndate_created = (2024, 5, 10)`nminimum_balance = 30"
$ws.Range("G12").Value = "Account Number: 350 Balance: $350.00\nMinimum Balance: $30.00 Account Type: Savings"

# Restore the original (custom) row heights for the rows that received multi-line
# text, since Excel auto-grows row height to fit wrapped/line-broken content.
$ws.Rows.Item(7).RowHeight = 31.2
$ws.Rows.Item(8).RowHeight = 31.2
$ws.Rows.Item(9).RowHeight = 31.2
$ws.Rows.Item(10).RowHeight = 31.2
$ws.Rows.Item(11).RowHeight = 31.2
$ws.Rows.Item(12).RowHeight = 49.5

# Leave the active cell on E12, matching where the final edits were made.
$ws.Range("E12").Select()
